# Apply updated simulation results (380 kV case) to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.4303915378680472
$ws.Range("C2").Value = 0.04382327218806381
$ws.Range("D2").Value = 0.2774905871002744
$ws.Range("F2").Value = 1.582889080602584
$ws.Range("G2").Value = 0.002470908429547307
$ws.Range("J2").Value = 0.3507128815763565
$ws.Range("K2").Value = 0.3972748393527752
$ws.Range("O2").Value = 3.69510745402863
$ws.Range("B3").Value = 0.3904979307735346
$ws.Range("C3").Value = 0.03851467051987356
$ws.Range("D3").Value = 0.266707547815983
$ws.Range("F3").Value = 1.58536303286801
$ws.Range("G3").Value = 0.002473572851304841
$ws.Range("J3").Value = 0.3394726189720245
$ws.Range("K3").Value = 0.3547135306349105
$ws.Range("O3").Value = 3.718957360956409
$ws.Range("B4").Value = 0.3660766694497397
$ws.Range("C4").Value = 0.0352395364321012
$ws.Range("D4").Value = 0.2602052215266326
$ws.Range("F4").Value = 1.587771657943016
$ws.Range("G4").Value = 0.002475295733736321
$ws.Range("J4").Value = 0.3327769460786101
$ws.Range("K4").Value = 0.3285974827209941
$ws.Range("O4").Value = 3.735687377764336
$ws.Range("B5").Value = 0.3561438333010472
$ws.Range("C5").Value = 0.03390103695186042
$ws.Range("D5").Value = 0.2575853901618359
$ws.Range("F5").Value = 1.58897683364367
$ws.Range("G5").Value = 0.002476019743942404
$ws.Range("J5").Value = 0.3301001377858768
$ws.Range("K5").Value = 0.3179597534052903
$ws.Range("O5").Value = 3.743029245639747
$ws.Range("B6").Value = 0.3544956580888652
$ws.Range("C6").Value = 0.03367854947215676
$ws.Range("D6").Value = 0.2571521797890028
$ws.Range("F6").Value = 1.589190457874452
$ws.Range("G6").Value = 0.002476141291214983
$ws.Range("J6").Value = 0.3296587791670191
$ws.Range("K6").Value = 0.316193671236249
$ws.Range("O6").Value = 3.74428001155178
$ws.Range("B7").Value = 0.3659426340137202
$ws.Range("C7").Value = 0.03522150047196249
$ws.Range("D7").Value = 0.2601697682616475
$ws.Range("F7").Value = 1.587787005928043
$ws.Range("G7").Value = 0.00247530540908857
$ws.Range("J7").Value = 0.3327406362942469
$ws.Range("K7").Value = 0.3284539985207005
$ws.Range("O7").Value = 3.735784270514472
$ws.Range("B8").Value = 0.4166212532286409
$ws.Range("C8").Value = 0.04199615373907761
$ws.Range("D8").Value = 0.2737480681930151
$ws.Range("F8").Value = 1.583557422560631
$ws.Range("G8").Value = 0.002471809122196536
$ws.Range("J8").Value = 0.3467944768161431
$ws.Range("K8").Value = 0.3825965392031208
$ws.Range("O8").Value = 3.702897746263403
$ws.Range("B9").Value = 0.5165677555208674
$ws.Range("C9").Value = 0.0551545514991858
$ws.Range("D9").Value = 0.301311662641325
$ws.Range("F9").Value = 1.582326316437289
$ws.Range("G9").Value = 0.002465639555038971
$ws.Range("J9").Value = 0.375991410857381
$ws.Range("K9").Value = 0.4888839106023681
$ws.Range("O9").Value = 3.654975425647706
$ws.Range("B10").Value = 0.5903270041803523
$ws.Range("C10").Value = 0.0647421791991718
$ws.Range("D10").Value = 0.322131055108656
$ws.Range("F10").Value = 1.58573649977383
$ws.Range("G10").Value = 0.0024615211511131
$ws.Range("J10").Value = 0.3984484545795652
$ws.Range("K10").Value = 0.5670251207922661
$ws.Range("O10").Value = 3.629892136702722
$ws.Range("B11").Value = 0.623950364150943
$ws.Range("C11").Value = 0.06908601186584917
$ws.Range("D11").Value = 0.3317253763217991
$ws.Range("F11").Value = 1.588226854490443
$ws.Range("G11").Value = 0.002459736663038641
$ws.Range("J11").Value = 0.4088850513699072
$ws.Range("K11").Value = 0.6025816436512343
$ws.Range("O11").Value = 3.620685292605629
$ws.Range("B12").Value = 0.6366922914097017
$ws.Range("C12").Value = 0.07072831240029132
$ws.Range("D12").Value = 0.3353761660865757
$ws.Range("F12").Value = 1.589305042919477
$ws.Range("G12").Value = 0.002459073653194891
$ws.Range("J12").Value = 0.4128689563800094
$ws.Range("K12").Value = 0.6160469253196368
$ws.Range("O12").Value = 3.617516208844762
$ws.Range("B13").Value = 0.6339476762161098
$ws.Range("C13").Value = 0.07037473114689874
$ws.Range("D13").Value = 0.3345891211753553
$ws.Range("F13").Value = 1.589066823359303
$ws.Range("G13").Value = 0.002459215878744256
$ws.Range("J13").Value = 0.4120095364683038
$ws.Range("K13").Value = 0.6131469099392177
$ws.Range("O13").Value = 3.618184606208445
$ws.Range("B14").Value = 0.6249984631796224
$ws.Range("C14").Value = 0.06922117777699555
$ws.Range("D14").Value = 0.3320253764176186
$ws.Range("F14").Value = 1.588312848553542
$ws.Range("G14").Value = 0.002459681861710669
$ws.Range("J14").Value = 0.409212172209692
$ws.Range("K14").Value = 0.6036894278825002
$ws.Range("O14").Value = 3.620418207392646
$ws.Range("B15").Value = 0.6195180311113972
$ws.Range("C15").Value = 0.06851425016920132
$ws.Range("D15").Value = 0.3304573007184217
$ws.Range("F15").Value = 1.587868620407534
$ws.Range("G15").Value = 0.002459968946886615
$ws.Range("J15").Value = 0.407502847461501
$ws.Range("K15").Value = 0.5978965333607107
$ws.Range("O15").Value = 3.621827694469744
$ws.Range("B16").Value = 0.5881309971373696
$ws.Range("C16").Value = 0.06445793797153954
$ws.Range("D16").Value = 0.3215065174245808
$ws.Range("F16").Value = 1.585592654772853
$ws.Range("G16").Value = 0.002461639557280065
$ws.Range("J16").Value = 0.3977708448564954
$ws.Range("K16").Value = 0.5647015624556957
$ws.Range("O16").Value = 3.63053820869203
$ws.Range("B17").Value = 0.5688935726000466
$ws.Range("C17").Value = 0.06196495002636482
$ws.Range("D17").Value = 0.3160470362657009
$ws.Range("F17").Value = 1.584437015088582
$ws.Range("G17").Value = 0.002462687173262449
$ws.Range("J17").Value = 0.3918571544514435
$ws.Range("K17").Value = 0.5443395808058256
$ws.Range("O17").Value = 3.636446569440977
$ws.Range("B18").Value = 0.5578353282814703
$ws.Range("C18").Value = 0.06052939486605169
$ws.Range("D18").Value = 0.3129185167450146
$ws.Range("F18").Value = 1.583860702427629
$ws.Range("G18").Value = 0.002463298114566803
$ws.Range("J18").Value = 0.3884765338542593
$ws.Range("K18").Value = 0.5326288553893335
$ws.Range("O18").Value = 3.640052294399311
$ws.Range("B19").Value = 0.5540923498073198
$ws.Range("C19").Value = 0.06004305884940209
$ws.Range("D19").Value = 0.3118612550687203
$ws.Range("F19").Value = 1.583680748936388
$ws.Range("G19").Value = 0.002463506409850791
$ws.Range("J19").Value = 0.3873354810172316
$ws.Range("K19").Value = 0.5286639859963884
$ws.Range("O19").Value = 3.641308736608323
$ws.Range("B20").Value = 0.5709407488968168
$ws.Range("C20").Value = 0.06223050465840174
$ws.Range("D20").Value = 0.3166270046525028
$ws.Range("F20").Value = 1.584550887020399
$ws.Range("G20").Value = 0.002462574785953131
$ws.Range("J20").Value = 0.3924845261497438
$ws.Range("K20").Value = 0.5465070547406583
$ws.Range("O20").Value = 3.635796146680008
$ws.Range("B21").Value = 0.6276268108745739
$ws.Range("C21").Value = 0.0695600759572983
$ws.Range("D21").Value = 0.3327779329765974
$ws.Range("F21").Value = 1.588530640542317
$ws.Range("G21").Value = 0.002459544646065993
$ws.Range("J21").Value = 0.4100329627090673
$ws.Range("K21").Value = 0.6064673026850187
$ws.Range("O21").Value = 3.619753528149033
$ws.Range("B22").Value = 0.6647294887715418
$ws.Range("C22").Value = 0.07433509791410131
$ws.Range("D22").Value = 0.3434361921884488
$ws.Range("F22").Value = 1.5919194204333
$ws.Range("G22").Value = 0.002457638489210487
$ws.Range("J22").Value = 0.4216872111009025
$ws.Range("K22").Value = 0.6456591690889866
$ws.Range("O22").Value = 3.611118713071221
$ws.Range("B23").Value = 0.6449222383162407
$ws.Range("C23").Value = 0.07178800284142994
$ws.Range("D23").Value = 0.3377383240594725
$ws.Range("F23").Value = 1.590038646039204
$ws.Range("G23").Value = 0.002458649070110286
$ws.Range("J23").Value = 0.4154501467889986
$ws.Range("K23").Value = 0.6247415277477728
$ws.Range("O23").Value = 3.615557851600585
$ws.Range("B24").Value = 0.5700152154448972
$ws.Range("C24").Value = 0.06211045457631315
$ws.Range("D24").Value = 0.3163647691365838
$ws.Range("F24").Value = 1.58449913116084
$ws.Range("G24").Value = 0.00246262556931872
$ws.Range("J24").Value = 0.3922008314850416
$ws.Range("K24").Value = 0.5455271532838992
$ws.Range("O24").Value = 3.636089552028579
$ws.Range("B25").Value = 0.4894707067339539
$ws.Range("C25").Value = 0.05160869157478487
$ws.Range("D25").Value = 0.2937549765996437
$ws.Range("F25").Value = 1.581902297716525
$ws.Range("G25").Value = 0.002467235514459021
$ws.Range("J25").Value = 0.3679167592233483
$ws.Range("K25").Value = 0.4601199913982725
$ws.Range("O25").Value = 3.666163385015665
$wb.Save()
